$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in column H, copying the formatting of the
# existing header cell (G1) so it matches the other header cells.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# Populate the "Save" column values for each data row
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
